$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.766.06'
$ws.Range("E2").Value = '  +2.70%  '

$ws.Range("D3").Value = '2.231.34'
$ws.Range("E3").Value = '  +0.85%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.89%  '

$ws.Range("E6").Value = '  -1.67%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.40'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.71%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.405'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.23%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.15'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.57%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0903'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.14%  '

$ws.Range("E12").Value = '  -0.19%  '

$ws.Range("D13").Value = '2.560.96'
$ws.Range("E13").Value = '  +0.93%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.75'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.81%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.00'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.21%  '

$ws.Range("E16").Value = '  -2.73%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.61'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.05%  '

$ws.Range("D18").Value = '2.240.81'
$ws.Range("E18").Value = '  +1.45%  '

$ws.Range("D19").Value = '41.729.38'
$ws.Range("E19").Value = '  +3.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.31'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.16%  '

$ws.Range("E21").Value = '  -0.44%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.11'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.34%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '248.73'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.50%  '

$ws.Range("E24").Value = '  -0.02%  '

$ws.Range("E25").Value = '  -0.05%  '

$ws.Range("E26").Value = '  +0.20%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.75'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.39%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '169.50'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.13%  '

$ws.Range("E29").Value = '  +0.80%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.94'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.23%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.42'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.86%  '

$ws.Range("E32").Value = '  -5.91%  '

$ws.Range("E33").Value = '  -1.09%  '

$ws.Range("E34").Value = '  +6.10%  '

$ws.Range("E35").Value = '  +0.06%  '

$ws.Range("E36").Value = '  +2.49%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.58'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.39%  '

$ws.Range("E38").Value = '  -4.76%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.39'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.75%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.00%  '

$ws.Range("B41").Value = 'TerraClassic'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.000237'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +13.83%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0240'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.41%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.52'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.14%  '

$ws.Range("E44").Value = '  -2.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '98.67'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.11%  '

$ws.Range("E46").Value = '  +2.72%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.38'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.90%  '

$ws.Range("D48").Value = '1.468.41'
$ws.Range("E48").Value = '  -3.56%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.59'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.87%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.78'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.33%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.28'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +11.35%  '
